# Add 2022-Q3 data
#
# 1. Insert a new worksheet named "2022-Q3" right before the existing
#    "2022-Q2" sheet (all sheets from 2022-Q2 onward keep their own
#    content unchanged, they just shift one tab to the right).
# 2. Populate the new "2022-Q3" sheet with its fund-holding table.
# 3. Insert a new row 2 in the "总计" (totals) sheet summarizing the
#    2022-Q3 quarter, and keep the running index in column A sequential.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" worksheet before "2022-Q2"
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($anchor)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# Step 2: fill in the "2022-Q3" fund-holding table
# ---------------------------------------------------------------------

# Header row (bold, bordered, centered - matches the other quarter sheets)
$header = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $header.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)   # starts at column B
    $cell.Value = $header[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
}

# Data rows - B..G are stored as text in the source data (fund codes,
# percentages, etc. are all kept as literal strings), H is numeric.
$rows = @(
    @{ idx = 0; code = "161914"; name = "万家创业板2年定期开放混合A"; scale = "7.30"; pos = "99.15"; ratio = "6.54"; value = "0.4774"; rank = 7 },
    @{ idx = 1; code = "161915"; name = "万家创业板2年定期开放混合C"; scale = "1.43"; pos = "99.15"; ratio = "6.54"; value = "0.0935"; rank = 7 }
)

$r = 2
foreach ($row in $rows) {
    $aCell = $q3.Cells.Item($r, 1)
    $aCell.Value = $row.idx
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $q3.Cells.Item($r, 2).Value = "'" + $row.code
    $q3.Cells.Item($r, 3).Value = $row.name
    $q3.Cells.Item($r, 4).Value = "'" + $row.scale
    $q3.Cells.Item($r, 5).Value = "'" + $row.pos
    $q3.Cells.Item($r, 6).Value = "'" + $row.ratio
    $q3.Cells.Item($r, 7).Value = "'" + $row.value
    $q3.Cells.Item($r, 8).Value = $row.rank

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 3: add the 2022-Q3 summary row to the "总计" sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Clear any formatting the row-insert may have pulled down from the
# header row, then restore it cell by cell to match the sheet's layout.
$total.Range("A2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.57

# Re-sequence the running index in column A for every row below, since
# they all shifted down by one position.
$idx = 1
for ($row = 3; $row -le 9; $row++) {
    $total.Cells.Item($row, 1).Value = $idx
    $idx = $idx + 1
}
